$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 15:52"

# Row 9 - Reino Unido (United Kingdom)
$ws.Range("B9").Value = 88621
$ws.Range("C9").Value = 4342
$ws.Range("E9").Value = 76948
$ws.Range("G9").Value = 717
$ws.Range("H9").Value = 11329

# Row 42 - Emiratos Arabes Unidos (United Arab Emirates)
$ws.Range("D42").Value = 852
$ws.Range("E42").Value = 3246
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 25

# Row 43 - Serbia
$ws.Range("B43").Value = 4054
$ws.Range("C43").Value = 424
$ws.Range("E43").Value = 3569
$ws.Range("F43").Value = 138
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 85

# Row 75 - Kazajistan (Kazakhstan)
$ws.Range("D75").Value = 138
$ws.Range("E75").Value = 829
